# refactor: split scraper logic into feature-based modules.
#
# The scraper's feature sections were re-ordered/split on the page being
# scraped, so the "taille" (size) filter block moved from section[5] to
# section[3] in the DOM. Update the xpath selectors that were anchored to
# that section accordingly: taille_label, douze_pouces and
# number_products (the latter is nested under the same section).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldPath = "section[1]/section[5]"
$newPath = "section[1]/section[3]"

$selectorNames = @("taille_label", "douze_pouces", "number_products")

foreach ($name in $selectorNames) {
    $nameCell = $ws.Columns.Item(1).Find($name)
    if ($nameCell -ne $null) {
        $xpathCell = $ws.Cells.Item($nameCell.Row, 2)
        $currentXpath = $xpathCell.Value()
        $xpathCell.Value = $currentXpath.Replace($oldPath, $newPath)
    }
}

# Minor row-height touch-up left over from the re-save (header row and the
# trailing "suivant" row), matching the workbook's own re-flow.
$ws.Rows.Item(1).RowHeight = 18.75
$ws.Rows.Item($ws.UsedRange.Rows.Count).RowHeight = 18.75
